$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.895.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.547.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.766.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.546.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.877.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("E22").Value = "  -2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.40%  "

$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.362.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.99%  "

$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("E36").Value = "  +4.43%  "

$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.76%  "

$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  +2.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("E46").Value = "  -2.59%  "

$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.681.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("E50").Value = "  +0.53%  "

$ws.Range("E51").Value = "  -0.81%  "
